# Update HOST column (C) for all test rows from "1PSEARCH" to "1PSEARCHV3"
# and update related QUERYSTRING (G) / API (D) values that reference the
# "wos" prefixed endpoints, per commit "Modified Test cases in to 1PSEARCHV3".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (HOST) : 1PSEARCH -> 1PSEARCHV3 for rows 2 through 25
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 3).Value = "1PSEARCHV3"
}

# Column G (QUERYSTRING) updates: add "wos::" prefix to the id lookup
$ws.Range("G3").Value = "?id=wos::(S1_TC_T1_hits.hits._id)"
$ws.Range("G4").Value = "?id=wos::(S1_TC_T1_hits.hits._id)&fields=category"

# Column D (API) updates: prefix /details paths with /wos
$ws.Range("D5").Value = "/wos/details/(S1_TC_T1_hits.hits._id)"
$ws.Range("D21").Value = "/wos/details/(S1_TC_T19_hits.hits._id)"
$ws.Range("D23").Value = "/wos/details/(S1_TC_T21_hits.hits._id)"
$ws.Range("D25").Value = "/wos/details/(S1_TC_T23_hits.hits._id)"
